$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 15 (TESAT_DATA / SDFSDFDSFSDFDS) - shifts subsequent rows up
$ws.Rows.Item(15).Delete()

# Update the selection to match target state
$ws.Range("B21").Select()
